$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 164
$ws.Range("F5").Value = 463
$ws.Range("F6").Value = 825
$ws.Range("F7").Value = 250
$ws.Range("F8").Value = 1209
$ws.Range("F9").Value = 347
$ws.Range("F12").Value = 693
$ws.Range("F13").Value = 187
$ws.Range("F14").Value = 514
$ws.Range("F18").Value = 2935
$ws.Range("F24").Value = 231
$ws.Range("F26").Value = 5297
$ws.Range("F28").Value = 987
$ws.Range("F29").Value = 23
$ws.Range("F31").Value = 316
$ws.Range("F32").Value = 1101
$ws.Range("F35").Value = 289

# Sheet 2: 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1125
$ws.Range("F10").Value = 31
$ws.Range("F24").Value = 319
$ws.Range("F25").Value = 278
$ws.Range("F26").Value = 3943
$ws.Range("F27").Value = 2
$ws.Range("F31").Value = 52

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 2466
$ws.Range("F6").Value = 1046
$ws.Range("F9").Value = 1329
$ws.Range("F10").Value = 364
$ws.Range("F11").Value = 102

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 2466
$ws.Range("F6").Value = 1046
$ws.Range("F7").Value = 1329
$ws.Range("F8").Value = 364
$ws.Range("F9").Value = 102
$ws.Range("F10").Value = 164
$ws.Range("F11").Value = 463
$ws.Range("F12").Value = 825
$ws.Range("F13").Value = 250
$ws.Range("F14").Value = 1209
$ws.Range("F15").Value = 347
$ws.Range("F17").Value = 693
$ws.Range("F18").Value = 1125
$ws.Range("F19").Value = 1125
$ws.Range("F20").Value = 187
$ws.Range("F21").Value = 514
$ws.Range("F23").Value = 2935
$ws.Range("F27").Value = 31
$ws.Range("F28").Value = 231
$ws.Range("F29").Value = 5297
$ws.Range("F31").Value = 987
$ws.Range("F34").Value = 23
$ws.Range("F37").Value = 316
$ws.Range("F43").Value = 319
$ws.Range("F44").Value = 319
$ws.Range("F45").Value = 278
$ws.Range("F46").Value = 1101
$ws.Range("F48").Value = 52
$ws.Range("F51").Value = 289
